$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 41 (shifts existing rows 41..150 down to 42..151)
$ws.Rows.Item(41).Insert()

# Populate the newly inserted row 41 with the new weekly data point
$ws.Cells.Item(41, 1).Value = 4
$ws.Cells.Item(41, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(41, 3).Value = "Los Lagos"
$ws.Cells.Item(41, 4).Value = 44497
$ws.Cells.Item(41, 5).Value = 10
$ws.Cells.Item(41, 6).Value = "Fruta"
$ws.Cells.Item(41, 7).Value = 100108
$ws.Cells.Item(41, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(41, 9).Value = 100108005
$ws.Cells.Item(41, 10).Value = "Piña"
$ws.Cells.Item(41, 11).Value = "Caramelo"
$ws.Cells.Item(41, 12).Value = "Segunda"
$ws.Cells.Item(41, 13).Value = 80
$ws.Cells.Item(41, 14).Value = 23000
$ws.Cells.Item(41, 15).Value = 24000
$ws.Cells.Item(41, 16).Value = 23500
$ws.Cells.Item(41, 17).Value = "`$/caja 14 unidades"
$ws.Cells.Item(41, 18).Value = "Ecuador"
$ws.Cells.Item(41, 19).Value = 1679
$ws.Cells.Item(41, 20).Value = 14
